$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.686.17"
$ws.Range("E2").Value = "  -0.84%  "
$ws.Range("D3").Value = "1.583.44"
$ws.Range("E3").Value = "  -3.17%  "
$ws.Range("E4").Value = "  +0.34%  "
$ws.Range("D5").Value = "'206.25"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.55%  "
$ws.Range("D6").Value = "'0.506"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.62%  "
$ws.Range("E7").Value = "  +0.35%  "
$ws.Range("D8").Value = "'22.23"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -5.07%  "
$ws.Range("E9").Value = "  -1.71%  "
$ws.Range("D10").Value = "'0.0591"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -3.20%  "
$ws.Range("D11").Value = "'0.0867"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.78%  "
$ws.Range("D12").Value = "1.806.56"
$ws.Range("E12").Value = "  -3.27%  "
$ws.Range("D13").Value = "1.588.14"
$ws.Range("E13").Value = "  -3.31%  "
$ws.Range("D14").Value = "'3.85"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -4.12%  "
$ws.Range("D15").Value = "'0.530"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -6.05%  "
$ws.Range("D16").Value = "27.638.82"
$ws.Range("E16").Value = "  -1.06%  "
$ws.Range("D17").Value = "'63.21"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -3.21%  "
$ws.Range("D18").Value = "'219.82"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -4.38%  "
$ws.Range("D19").Value = "0.0₃0691"
$ws.Range("E19").Value = "  -3.96%  "
$ws.Range("D20").Value = "'7.31"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -6.75%  "
$ws.Range("D22").Value = "'4.13"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -5.50%  "
$ws.Range("D23").Value = "'9.49"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -6.59%  "
$ws.Range("D24").Value = "'1.97"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -5.46%  "
$ws.Range("D25").Value = "'153.59"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.74%  "
$ws.Range("E26").Value = "  +0.37%  "
$ws.Range("D27").Value = "'6.75"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -3.22%  "
$ws.Range("D28").Value = "'15.13"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -2.80%  "
$ws.Range("E29").Value = "  -4.07%  "
$ws.Range("E30").Value = "  -2.80%  "
$ws.Range("D31").Value = "'0.0464"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -3.68%  "
$ws.Range("D32").Value = "'3.22"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -5.66%  "
$ws.Range("D33").Value = "1.386.18"
$ws.Range("E33").Value = "  -1.15%  "
$ws.Range("E34").Value = "  -5.60%  "
$ws.Range("D35").Value = "'1.52"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -5.59%  "
$ws.Range("D36").Value = "'0.971"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -4.12%  "
$ws.Range("E37").Value = "  -0.80%  "
$ws.Range("E38").Value = "  -3.54%  "
$ws.Range("D39").Value = "'0.539"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -3.63%  "
$ws.Range("D40").Value = "'0.818"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -4.08%  "
$ws.Range("E41").Value = "  +0.38%  "
$ws.Range("D42").Value = "'0.977"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -2.45%  "
$ws.Range("D43").Value = "'2.17"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +1.38%  "
$ws.Range("E44").Value = "  -4.74%  "
$ws.Range("D45").Value = "'63.46"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -4.01%  "
$ws.Range("D46").Value = "'5.21"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -4.57%  "
$ws.Range("D47").Value = "1.718.93"
$ws.Range("E47").Value = "  -3.29%  "
$ws.Range("D48").Value = "'87.97"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.66%  "
$ws.Range("E49").Value = "  -1.94%  "
$ws.Range("D50").Value = "'0.0973"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -5.02%  "
$ws.Range("E51").Value = "  -1.29%  "
